$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("G2").Value = 8.142376000000001
$ws.Range("H2").Value = 24.427128
$ws.Range("I2").Value = 0.1741313933276368
$ws.Range("J2").Value = 0.1741313933276368
$ws.Range("M2").Value = 0.4890553333333333
$ws.Range("N2").Value = 1.467166
$ws.Range("O2").Value = 0.9644476581758422
$ws.Range("P2").Value = 0.9644476581758422
$ws.Range("Q2").Value = 3.982072408805334
$ws.Range("R2").Value = 35.838651679248
$ws.Range("S2").Value = 0.1679406145097358
$ws.Range("T2").Value = 0.1679406145097358

# Row 3
$ws.Range("G3").Value = 8.142376000000001
$ws.Range("H3").Value = 24.427128
$ws.Range("I3").Value = 0.1741313933276368
$ws.Range("J3").Value = 0.1741313933276368
$ws.Range("O3").Value = 0.03555234182415776
$ws.Range("P3").Value = 0.03555234182415776
$ws.Range("Q3").Value = 0.146790754528
$ws.Range("R3").Value = 1.321116790752
$ws.Range("S3").Value = 0.006190778817901007
$ws.Range("T3").Value = 0.006190778817901007

# Row 4
$ws.Range("I4").Value = 0.5205382400466131
$ws.Range("J4").Value = 0.5205382400466131
$ws.Range("M4").Value = 0.4890553333333333
$ws.Range("N4").Value = 1.467166
$ws.Range("O4").Value = 0.9644476581758422
$ws.Range("P4").Value = 0.9644476581758422
$ws.Range("Q4").Value = 11.90377521138645
$ws.Range("R4").Value = 107.133976902478
$ws.Range("S4").Value = 0.5020318866039304
$ws.Range("T4").Value = 0.5020318866039304

# Row 5
$ws.Range("I5").Value = 0.5205382400466131
$ws.Range("J5").Value = 0.5205382400466131
$ws.Range("O5").Value = 0.03555234182415776
$ws.Range("P5").Value = 0.03555234182415776
$ws.Range("S5").Value = 0.01850635344268267
$ws.Range("T5").Value = 0.01850635344268267

# Row 6
$ws.Range("I6").Value = 0.3053303666257501
$ws.Range("J6").Value = 0.3053303666257501
$ws.Range("M6").Value = 0.4890553333333333
$ws.Range("N6").Value = 1.467166
$ws.Range("O6").Value = 0.9644476581758422
$ws.Range("P6").Value = 0.9644476581758422
$ws.Range("Q6").Value = 6.982357432948
$ws.Range("R6").Value = 62.841216896532
$ws.Range("S6").Value = 0.294475157062176
$ws.Range("T6").Value = 0.294475157062176

# Row 7
$ws.Range("I7").Value = 0.3053303666257501
$ws.Range("J7").Value = 0.3053303666257501
$ws.Range("O7").Value = 0.03555234182415776
$ws.Range("P7").Value = 0.03555234182415776
$ws.Range("S7").Value = 0.01085520956357408
$ws.Range("T7").Value = 0.01085520956357408
